# Scheduled runner update: refresh market price snapshot columns (H:N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 345.4
$ws.Range("I2").Value = 398
$ws.Range("J2").Value = 310.33334
$ws.Range("K2").Value = 398
$ws.Range("L2").Value = 310.33334
$ws.Range("M2").Value = -285
$ws.Range("N2").Value = -536.33334

$ws.Range("H9").Value = 125002500
$ws.Range("I9").Value = 125002500
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 125002500
$ws.Range("L9").Value = 0
$ws.Range("N9").Value = -125002331
$ws.Range("M9").ClearContents()

$ws.Range("H21").Value = 12499.75
$ws.Range("I21").Value = 12499.75
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 12499.75
$ws.Range("L21").Value = 0
$ws.Range("N21").Value = -12031.75
$ws.Range("M21").ClearContents()

$ws.Range("H23").Value = 12499.75
$ws.Range("I23").Value = 12499.75
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 12499.75
$ws.Range("L23").Value = 0
$ws.Range("N23").Value = -12265.75
$ws.Range("M23").ClearContents()

$ws.Range("H38").Value = 424.1111
$ws.Range("I38").Value = 142.5
$ws.Range("J38").Value = 987.3333
$ws.Range("K38").Value = 427.5
$ws.Range("L38").Value = 2961.9999
$ws.Range("M38").Value = -55.5
$ws.Range("N38").Value = -3705.9999

$ws.Range("H43").Value = 799.2
$ws.Range("I43").Value = 785.2857
$ws.Range("J43").Value = 831.6667
$ws.Range("K43").Value = 785.2857
$ws.Range("L43").Value = 831.6667
$ws.Range("M43").Value = -716.2857
$ws.Range("N43").Value = -969.6667

$ws.Range("H51").Value = 5796.316
$ws.Range("I51").Value = 9000.5
$ws.Range("J51").Value = 5419.353
$ws.Range("K51").Value = 9000.5
$ws.Range("L51").Value = 5419.353
$ws.Range("M51").Value = -8516.5
$ws.Range("N51").Value = -6387.353

$ws.Range("H70").Value = 930333.4399999999
$ws.Range("I70").Value = 2042615
$ws.Range("J70").Value = 3432.1667
$ws.Range("K70").Value = 6127845
$ws.Range("L70").Value = 10296.5001
$ws.Range("M70").Value = -6127575
$ws.Range("N70").Value = -10836.5001

$ws.Range("H73").Value = 930333.4399999999
$ws.Range("I73").Value = 2042615
$ws.Range("J73").Value = 3432.1667
$ws.Range("K73").Value = 6127845
$ws.Range("L73").Value = 10296.5001
$ws.Range("M73").Value = -6126909
$ws.Range("N73").Value = -12168.5001

$ws.Range("H86").Value = 20106948
$ws.Range("I86").Value = 4895.8
$ws.Range("J86").Value = 40209000
$ws.Range("K86").Value = 4895.8
$ws.Range("L86").Value = 40209000
$ws.Range("M86").Value = -3772.8
$ws.Range("N86").Value = -40211246

$ws.Range("H89").Value = 20106948
$ws.Range("I89").Value = 4895.8
$ws.Range("J89").Value = 40209000
$ws.Range("K89").Value = 24479
$ws.Range("L89").Value = 201045000
$ws.Range("M89").Value = -18863
$ws.Range("N89").Value = -201056232

$ws.Range("H98").Value = 538.1579
$ws.Range("I98").Value = 555.8823
$ws.Range("J98").Value = 387.5
$ws.Range("K98").Value = 555.8823
$ws.Range("L98").Value = 387.5
$ws.Range("M98").Value = 942.1177
$ws.Range("N98").Value = -3383.5

$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("N105").Value = 0
$ws.Range("L105").ClearContents()

$ws.Range("H111").Value = 110031.5
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 110031.5
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 330094.5
$ws.Range("N111").Value = -336228.5
$ws.Range("M111").ClearContents()

$ws.Range("H122").Value = 538.1579
$ws.Range("I122").Value = 555.8823
$ws.Range("J122").Value = 387.5
$ws.Range("K122").Value = 1667.6469
$ws.Range("L122").Value = 1162.5
$ws.Range("M122").Value = 782.3531
$ws.Range("N122").Value = -6062.5

$ws.Range("H125").Value = 1246.2222
$ws.Range("I125").Value = 934.8333
$ws.Range("J125").Value = 1869
$ws.Range("K125").Value = 8413.4997
$ws.Range("L125").Value = 16821
$ws.Range("M125").Value = -5953.4997
$ws.Range("N125").Value = -21741

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 17494.8
$ws.Range("I31").Value = 9368.5
$ws.Range("J31").Value = 50000
$ws.Range("K31").Value = 9368.5
$ws.Range("L31").Value = 50000
$ws.Range("M31").Value = -9074.5
$ws.Range("N31").Value = -50588

$ws.Range("H32").Value = 2112.621
$ws.Range("I32").Value = 1298.0892
$ws.Range("J32").Value = 6674
$ws.Range("K32").Value = 1298.0892
$ws.Range("L32").Value = 6674
$ws.Range("M32").Value = -1011.0892
$ws.Range("N32").Value = -7248

$ws.Range("H63").Value = 85719650
$ws.Range("I63").Value = 125006056
$ws.Range("J63").Value = 33337766
$ws.Range("K63").Value = 125006056
$ws.Range("L63").Value = 33337766
$ws.Range("M63").Value = -125005370
$ws.Range("N63").Value = -33339138

$ws.Range("H66").Value = 85719650
$ws.Range("I66").Value = 125006056
$ws.Range("J66").Value = 33337766
$ws.Range("K66").Value = 625030280
$ws.Range("L66").Value = 166688830
$ws.Range("M66").Value = -625026848
$ws.Range("N66").Value = -166695694

$ws.Range("H74").Value = 50004236
$ws.Range("I74").Value = 83336570
$ws.Range("J74").Value = 5744.5
$ws.Range("K74").Value = 83336570
$ws.Range("L74").Value = 5744.5
$ws.Range("M74").Value = -83335696
$ws.Range("N74").Value = -7492.5

$ws.Range("H77").Value = 50004236
$ws.Range("I77").Value = 83336570
$ws.Range("J77").Value = 5744.5
$ws.Range("K77").Value = 416682850
$ws.Range("L77").Value = 28722.5
$ws.Range("M77").Value = -416678482
$ws.Range("N77").Value = -37458.5

$ws.Range("H102").Value = 2121363.2
$ws.Range("I102").Value = 2218655.2
$ws.Range("J102").Value = 126874.5
$ws.Range("K102").Value = 2218655.2
$ws.Range("L102").Value = 126874.5
$ws.Range("M102").Value = -2217033.2
$ws.Range("N102").Value = -130118.5

$ws.Range("H122").Value = 17546582
$ws.Range("I122").Value = 27779340
$ws.Range("J122").Value = 4714.143
$ws.Range("K122").Value = 83338020
$ws.Range("L122").Value = 14142.429
$ws.Range("M122").Value = -83335570
$ws.Range("N122").Value = -19042.429

$ws.Range("H132").Value = 22224026
$ws.Range("I132").Value = 26317482
$ws.Range("J132").Value = 2398
$ws.Range("K132").Value = 78952446
$ws.Range("L132").Value = 7194
$ws.Range("M132").Value = -78949916
$ws.Range("N132").Value = -12254

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H74").Value = 35852.145
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 35852.145
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 35852.145
$ws.Range("N74").Value = -37724.145

$ws.Range("H77").Value = 35852.145
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 35852.145
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 107556.435
$ws.Range("N77").Value = -116916.435

$ws.Range("H94").Value = 4148.9585
$ws.Range("I94").Value = 5245.7334
$ws.Range("J94").Value = 2321
$ws.Range("K94").Value = 5245.7334
$ws.Range("L94").Value = 2321
$ws.Range("M94").Value = -4794.7334
$ws.Range("N94").Value = -3223

$ws.Range("H105").Value = 1857.4
$ws.Range("I105").Value = 1757.9445
$ws.Range("J105").Value = 2113.1428
$ws.Range("K105").Value = 1757.9445
$ws.Range("L105").Value = 2113.1428
$ws.Range("M105").Value = -10.94450000000006
$ws.Range("N105").Value = -5607.1428

$ws.Range("H134").Value = 1322.5264
$ws.Range("I134").Value = 1118.2222
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 3354.6666
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -819.6665999999996
$ws.Range("N134").Value = -20070

$ws.Range("H137").Value = 54999.5
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 54999.5
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 54999.5
$ws.Range("N137").Value = -65199.5

$ws.Range("H138").Value = 59236
$ws.Range("I138").Value = 40709
$ws.Range("J138").Value = 68499.5
$ws.Range("K138").Value = 40709
$ws.Range("L138").Value = 68499.5
$ws.Range("M138").Value = -35569
$ws.Range("N138").Value = -78779.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 100.78571
$ws.Range("I7").Value = 109.166664
$ws.Range("J7").Value = 50.5
$ws.Range("K7").Value = 109.166664
$ws.Range("L7").Value = 50.5
$ws.Range("M7").Value = 3.833336000000003
$ws.Range("N7").Value = -276.5

$ws.Range("H31").Value = 3318.8306
$ws.Range("I31").Value = 1740.25
$ws.Range("J31").Value = 3566.451
$ws.Range("K31").Value = 1740.25
$ws.Range("L31").Value = 3566.451
$ws.Range("M31").Value = -1445.25
$ws.Range("N31").Value = -4156.451

$ws.Range("H34").Value = 3318.8306
$ws.Range("I34").Value = 1740.25
$ws.Range("J34").Value = 3566.451
$ws.Range("K34").Value = 1740.25
$ws.Range("L34").Value = 3566.451
$ws.Range("M34").Value = -1538.25
$ws.Range("N34").Value = -3970.451

$ws.Range("H58").Value = 2142.6667
$ws.Range("I58").Value = 2098
$ws.Range("J58").Value = 2500
$ws.Range("K58").Value = 2098
$ws.Range("L58").Value = 2500
$ws.Range("M58").Value = -1895
$ws.Range("N58").Value = -2906

$ws.Range("H94").Value = 1583.2142
$ws.Range("I94").Value = 1273
$ws.Range("J94").Value = 1707.3
$ws.Range("K94").Value = 1273
$ws.Range("L94").Value = 1707.3
$ws.Range("M94").Value = -822
$ws.Range("N94").Value = -2609.3

$ws.Range("H136").Value = 2142.6667
$ws.Range("I136").Value = 2098
$ws.Range("J136").Value = 2500
$ws.Range("K136").Value = 6294
$ws.Range("L136").Value = 7500
$ws.Range("M136").Value = -3744
$ws.Range("N136").Value = -12600

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 749.5
$ws.Range("I20").Value = 499
$ws.Range("J20").Value = 1000
$ws.Range("K20").Value = 1497
$ws.Range("L20").Value = 3000
$ws.Range("M20").Value = -1270
$ws.Range("N20").Value = -3454

$ws.Range("H26").Value = 528.3333
$ws.Range("I26").Value = 528.3333
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 1584.9999
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -1296.9999

$ws.Range("H31").Value = 2042.5
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 2042.5
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 6127.5
$ws.Range("N31").Value = -6703.5

$ws.Range("H34").Value = 1987.7142
$ws.Range("I34").Value = 797.3333
$ws.Range("J34").Value = 2186.111
$ws.Range("K34").Value = 2391.9999
$ws.Range("L34").Value = 6558.333
$ws.Range("M34").Value = -2307.9999
$ws.Range("N34").Value = -6726.333

$ws.Range("H86").Value = 771.125
$ws.Range("I86").Value = 748.75
$ws.Range("J86").Value = 793.5
$ws.Range("K86").Value = 2246.25
$ws.Range("L86").Value = 2380.5
$ws.Range("M86").Value = -1060.25
$ws.Range("N86").Value = -4752.5

$ws.Range("H89").Value = 771.125
$ws.Range("I89").Value = 748.75
$ws.Range("J89").Value = 793.5
$ws.Range("K89").Value = 6738.75
$ws.Range("L89").Value = 7141.5
$ws.Range("M89").Value = -810.75
$ws.Range("N89").Value = -18997.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4999.8335
$ws.Range("I80").Value = 2999.75
$ws.Range("J80").Value = 9000
$ws.Range("K80").Value = 2999.75
$ws.Range("L80").Value = 9000
$ws.Range("M80").Value = -2001.75
$ws.Range("N80").Value = -10996

$ws.Range("H83").Value = 4999.8335
$ws.Range("I83").Value = 2999.75
$ws.Range("J83").Value = 9000
$ws.Range("K83").Value = 14998.75
$ws.Range("L83").Value = 45000
$ws.Range("M83").Value = -10006.75
$ws.Range("N83").Value = -54984

$ws.Range("H99").Value = 1522
$ws.Range("I99").Value = 1522
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1522
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 724

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 6899.4736
$ws.Range("I132").Value = 3245.6667
$ws.Range("J132").Value = 10187.9
$ws.Range("K132").Value = 9737.000100000001
$ws.Range("L132").Value = 30563.7
$ws.Range("M132").Value = -7207.000100000001
$ws.Range("N132").Value = -35623.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5535.0527
$ws.Range("I132").Value = 5293.1333
$ws.Range("J132").Value = 6442.25
$ws.Range("K132").Value = 15879.3999
$ws.Range("L132").Value = 19326.75
$ws.Range("M132").Value = -13349.3999
$ws.Range("N132").Value = -24386.75
